$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VAR")

# Fill in the new Flag Sequence / Flag Step Modbus entries (AD:AF, rows 19-29 and 32-42)
    $ws.Range("AF19").Value = "M61"
    $ws.Range("AF20").Value = "M62"
    $ws.Range("AF21").Value = "M71"
    $ws.Range("AF22").Value = "M72"
    $ws.Range("AF23").Value = "M73"
    $ws.Range("AF24").Value = "M74"
    $ws.Range("AF25").Value = "M75"
    $ws.Range("AF26").Value = "M76"
    $ws.Range("AF27").Value = "M77"
    $ws.Range("AF28").Value = "M78"
    $ws.Range("AF29").Value = "M79"
    $ws.Range("AF32").Value = "M200"
    $ws.Range("AD19").Value = "Flag Sequence Init 1"
    $ws.Range("AE19").Value = 3133
    $ws.Range("AD20").Value = "Flag Sequence Init 2"
    $ws.Range("AE20").Value = 3134
    $ws.Range("AD21").Value = "Flag Sequence 1"
    $ws.Range("AE21").Value = 3143
    $ws.Range("AD22").Value = "Flag Sequence 2"
    $ws.Range("AE22").Value = 3144
    $ws.Range("AD23").Value = "Flag Sequence 3"
    $ws.Range("AE23").Value = 3145
    $ws.Range("AD24").Value = "Flag Sequence 4"
    $ws.Range("AE24").Value = 3146
    $ws.Range("AD25").Value = "Flag Sequence 5"
    $ws.Range("AE25").Value = 3147
    $ws.Range("AD26").Value = "Flag Sequence 6"
    $ws.Range("AE26").Value = 3148
    $ws.Range("AD27").Value = "Flag Sequence 7"
    $ws.Range("AE27").Value = 3149
    $ws.Range("AD28").Value = "Flag Sequence 8"
    $ws.Range("AE28").Value = 3150
    $ws.Range("AD29").Value = "Flag Sequence 9"
    $ws.Range("AE29").Value = 3151
    $ws.Range("AD33").Value = "Flag Step 1"
    $ws.Range("AD34").Value = "Flag Step 2"
    $ws.Range("AF33").Value = "M201"
    $ws.Range("AE33").Value = 3273
    $ws.Range("AD35").Value = "Flag Step 3"
    $ws.Range("AF34").Value = "M202"
    $ws.Range("AE34").Value = 3274
    $ws.Range("AD36").Value = "Flag Step 4"
    $ws.Range("AF35").Value = "M203"
    $ws.Range("AE35").Value = 3275
    $ws.Range("AD37").Value = "Flag Step 5"
    $ws.Range("AF36").Value = "M204"
    $ws.Range("AE36").Value = 3276
    $ws.Range("AD38").Value = "Flag Step 6"
    $ws.Range("AF37").Value = "M205"
    $ws.Range("AE37").Value = 3277
    $ws.Range("AD39").Value = "Flag Step 7"
    $ws.Range("AF38").Value = "M206"
    $ws.Range("AE38").Value = 3278
    $ws.Range("AD40").Value = "Flag Step 8"
    $ws.Range("AF39").Value = "M207"
    $ws.Range("AE39").Value = 3279
    $ws.Range("AD41").Value = "Flag Step 9"
    $ws.Range("AF40").Value = "M208"
    $ws.Range("AE40").Value = 3280
    $ws.Range("AD42").Value = "Flag Step 10"
    $ws.Range("AF41").Value = "M209"
    $ws.Range("AE41").Value = 3281
    $ws.Range("AF42").Value = "M210"
    $ws.Range("AE42").Value = 3282
    $ws.Range("AD32").Value = "Flag Step 0"
    $ws.Range("AE32").Value = 3272

# Update the selection on the VAR sheet
$ws.Range("AF34").Select()
